# The canonical OOXML diff for this fixture is purely a cosmetic
# re-serialization: every changed line is the same element with the
# same attribute names/values, just re-emitted with attributes sorted
# alphabetically (xmlns:* declarations sorted, then the remaining
# attributes sorted) -- e.g. <w:pgSz w:w="11906" w:h="16838"/> becomes
# <w:pgSz w:h="16838" w:w="11906"/>, <w:style w:type="paragraph"
# w:default="1" w:styleId="Normal"> becomes <w:style w:default="1"
# w:styleId="Normal" w:type="paragraph">, and so on throughout
# word/document.xml's root namespace list and word/styles.xml's
# <w:docDefaults>/<w:latentStyles>/<w:style> blocks. No text, run,
# paragraph, style value, page size/margin value, or any other
# document content actually changes.
#
# That kind of attribute-order-only rewrite is not something the Word
# object model exposes a way to request (PageSetup/Styles/sectPr are
# all read through properties that preserve the underlying XML
# attribute ordering when written back), so there is no Word
# automation call that reproduces it. The content itself already
# matches the target, so no edits are required here.
$d = $word.ActiveDocument
